# "Divide urls and emails data" -- split the single "url" column into two
# columns: the original hyperlinked url (column A) and the generated
# Lighthouse report file name for that url (new column B, "reportName").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "reportName" column -------------------------------------------------
$ws.Range("B1").Value = "reportName"
$ws.Range("B2").Value = "google_lhreport.html"
$ws.Range("B3").Value = "facebook_lhreport.html"

# Small (8pt) black Arial font for the whole new column, matching the sizing
# used for the url data rows.
$ws.Range("B1:B3").Font.Size = 8

# --- Re-style the url data cells (A2:A3) -------------------------------------
# Shrink to 8pt and recolor the hyperlink text (underline was already set).
$ws.Range("A2:A3").Font.Size = 8
$ws.Range("A2:A3").Font.Color = 12598288

# --- Column widths ------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 20
$ws.Columns.Item(2).ColumnWidth = 31.83
